# Auto-generated edit script: updates crypto price/volume values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.299.94'
Set-TextValue 'E2' '  +1.24%  '
Set-TextValue 'D3' '1.620.35'
Set-TextValue 'E3' '  +2.03%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '212.10'
Set-TextValue 'E5' '  +0.75%  '
Set-TextValue 'D7' '0.483'
Set-TextValue 'E7' '  +0.92%  '
Set-TextValue 'E8' '  +0.77%  '
Set-TextValue 'E9' '  +0.74%  '
Set-TextValue 'D10' '18.78'
Set-TextValue 'E10' '  +4.83%  '
Set-TextValue 'D11' '0.0815'
Set-TextValue 'E11' '  +0.89%  '
Set-TextValue 'D12' '1.845.43'
Set-TextValue 'E12' '  +1.97%  '
Set-TextValue 'D13' '1.624.05'
Set-TextValue 'E13' '  +2.34%  '
Set-TextValue 'D14' '4.01'
Set-TextValue 'E14' '  +0.67%  '
Set-TextValue 'D15' '0.518'
Set-TextValue 'E15' '  +1.61%  '
Set-TextValue 'D16' '26.303.58'
Set-TextValue 'E16' '  +1.31%  '
Set-TextValue 'D17' '62.21'
Set-TextValue 'E17' '  +3.56%  '
Set-TextValue 'E18' '  +0.81%  '
Set-TextValue 'E19' '  +0.00%  '
Set-TextValue 'D20' '201.68'
Set-TextValue 'E20' '  +1.03%  '
Set-TextValue 'E21' '  +1.66%  '
Set-TextValue 'E22' '  +1.69%  '
Set-TextValue 'E23' '  +1.25%  '
Set-TextValue 'E24' '  +0.42%  '
Set-TextValue 'D25' '144.69'
Set-TextValue 'E25' '  +1.43%  '
Set-TextValue 'E26' '  +0.03%  '
Set-TextValue 'E27' '  -1.27%  '
Set-TextValue 'D28' '15.16'
Set-TextValue 'E29' '  +1.67%  '
Set-TextValue 'D30' '0.0521'
Set-TextValue 'E30' '  +10.26%  '
Set-TextValue 'E31' '  +0.87%  '
Set-TextValue 'E32' '  +2.07%  '
Set-TextValue 'E33' '  -0.44%  '
Set-TextValue 'D34' '1.50'
Set-TextValue 'E34' '  +1.88%  '
Set-TextValue 'E35' '  +2.51%  '
Set-TextValue 'D36' '1.179.14'
Set-TextValue 'E36' '  +4.95%  '
Set-TextValue 'E37' '  +0.64%  '
Set-TextValue 'E38' '  +3.04%  '
Set-TextValue 'E39' '  -0.01%  '
Set-TextValue 'D40' '2.31'
Set-TextValue 'E40' '  +0.14%  '
Set-TextValue 'D41' '0.496'
Set-TextValue 'E41' '  +1.65%  '
Set-TextValue 'E42' '  +1.39%  '
Set-TextValue 'E43' '  +4.80%  '
Set-TextValue 'D44' '1.756.93'
Set-TextValue 'E44' '  +2.10%  '
Set-TextValue 'D45' '92.64'
Set-TextValue 'E45' '  +0.76%  '
Set-TextValue 'E46' '  +3.16%  '
Set-TextValue 'D47' '53.83'
Set-TextValue 'E47' '  +1.07%  '
Set-TextValue 'E48' '  +1.15%  '
Set-TextValue 'E49' '  +0.33%  '
Set-TextValue 'E50' '  -0.15%  '
Set-TextValue 'D51' '7.26'
Set-TextValue 'E51' '  +2.33%  '
